$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)   # "KPI - VISA"

# ---- Header row (A1:F1): bold + wrap text ----
$ws.Range("A1:F1").Font.Bold = $true
$ws.Range("A1:F1").WrapText = $true

# ---- Column widths ----
$ws.Columns.Item(2).ColumnWidth = 17.5    # -> stored width ~18.33 (target 18.34765625)
$ws.Columns.Item(3).ColumnWidth = 24      # -> stored width ~24.83 (target 24.84765625)
$ws.Columns.Item(4).ColumnWidth = 22.5    # -> stored width ~23.33 (target 23.296875)
$ws.Columns.Item(5).ColumnWidth = 13      # -> stored width ~13.83 (target 13.8984375)

# ---- Rows 2 & 3: columns A:D get "vertical top + wrap text" formatting ----
$ws.Range("A2:D3").WrapText = $true
$ws.Range("A2:D3").VerticalAlignment = -4160   # xlTop

$ws.Range("A2").Value2 = 1
$ws.Range("B2").Value2 = "To verify the records in database matches the result shown on web UI "
$ws.Range("A3").Value2 = 2
$ws.Range("B3").Value2 = "To verify the records in database matches the result shown on web UI "

# New test-case text (order matters for shared-string indices)
$ws.Range("C2").Value2 = "team = visa`nstart date = 2016 Oct 1`nend date = 2016 Oct 15"
$ws.Range("C3").Value2 = "team = visa`nstart date = 2016 Oct 15`nend date = 2016 Oct 31"
$ws.Range("D2").Value2 = "indonesian count = 0 `nnon-indonesian count = 0"
$ws.Range("D3").Value2 = "indonesian count = 0 `nnon-indonesian count = 2"

# ---- Row 4: B4:D4 get "vertical top + wrap text"; A4 gets wrap-text only ----
$ws.Range("B4:D4").WrapText = $true
$ws.Range("B4:D4").VerticalAlignment = -4160   # xlTop
$ws.Range("B4").Value2 = "To verify the records in database matches the result shown on web UI "
$ws.Range("C4").Value2 = "team = visa`nstart date = 2016 Nov 1`nend date = 2017 Dec 31"
$ws.Range("D4").Value2 = "indonesian count = 2 `nnon-indonesian count = 2"

$ws.Range("A4").WrapText = $true
$ws.Range("A4").Value2 = 3

# ---- Empty trailing cells E2:F4: wrap text only ----
$ws.Range("E2:F4").WrapText = $true

# ---- Row heights for the new data rows ----
$ws.Rows.Item(2).RowHeight = 62.4
$ws.Rows.Item(3).RowHeight = 62.4
$ws.Rows.Item(4).RowHeight = 62.4

# ---- Sheet4 selection state ----
$ws.Range("A1:F4").Select() | Out-Null

# ---- Switch the active tab from "Gender and Age" to "KPI - Medical" ----
$wsMedical = $wb.Worksheets.Item(5)   # "KPI - Medical"
$wsMedical.Activate()

Write-Host "edit applied"
